$d = $word.ActiveDocument

# Locate the anchor paragraph: "Semi-structured " (last paragraph of the Interview section)
$rng = $d.Content
$found = $rng.Find.Execute("Semi-structured ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find anchor paragraph Semi-structured" }
$anchorPara = $rng.Paragraphs(1)
$idx = $anchorPara.Index

# Insert the new "Report" section paragraphs after the anchor, preserving the
# existing trailing blank paragraph that follows it.
# --- paragraph 0 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# --- paragraph 1 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# --- paragraph 2 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Report: </w:t></w:r></w:p>')

# --- paragraph 3 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Aim: To critically reflect on the use of two key human geography methods (auto-ethnography and interviews), to analyse personal fieldwork material to discuss how your field site might be understood through geographical ideas of ‘Senses of Home’ </w:t></w:r></w:p>')

# --- paragraph 4 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# --- paragraph 5 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Structure: </w:t></w:r></w:p>')

# --- paragraph 6 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Title: Investigating senses of ‘home’</w:t></w:r></w:p>')

# --- paragraph 7 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Research Question: What critical understandings of ‘home’ can be gained through ethnographic research? </w:t></w:r></w:p>')

# --- paragraph 8 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">(200 Words) Introduction: Introduce field work, designated locations, particular geographical aspects of the site </w:t></w:r></w:p>')

# --- paragraph 9 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">Why is it relevant to study ‘Senses of home’ What can be learned from this study? </w:t></w:r></w:p>')

# --- paragraph 10 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Where did the fieldwork take place? Describe the area and spaces focused on</w:t></w:r></w:p>')

# --- paragraph 11 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">How has the report been focused? Particular approaches? Particular focuses? </w:t></w:r></w:p>')

# --- paragraph 12 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>(400 Words) Auto-ethnography reflection: Critically reflect on uses and challenges of auto-ethnographic observation</w:t></w:r></w:p>')

# --- paragraph 13 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Reflect on strengths and weaknesses of approach to method</w:t></w:r></w:p>')

# --- paragraph 14 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">Self-reflexivity </w:t></w:r></w:p>')

# --- paragraph 15 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Engage with literature</w:t></w:r></w:p>')

# --- paragraph 16 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">(400 Words) Interview reflection: </w:t></w:r><w:r><w:t>Critically reflect on uses and challenges</w:t></w:r><w:r><w:t xml:space="preserve"> of interviews</w:t></w:r></w:p>')

# --- paragraph 17 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Reflect on strengths and weaknesses of approach to method</w:t></w:r></w:p>')

# --- paragraph 18 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">Engage with literature </w:t></w:r></w:p>')

# --- paragraph 19 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">(800 Words) Fieldwork findings: Critically analyse the particular ways ‘Senses of home’ matter </w:t></w:r></w:p>')

# --- paragraph 20 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">Consult geographical ideas of home </w:t></w:r></w:p>')

# --- paragraph 21 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Include quotations</w:t></w:r></w:p>')

# --- paragraph 22 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Analyse ethnographic observations and interview</w:t></w:r></w:p>')

# --- paragraph 23 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Include visuals/notes</w:t></w:r></w:p>')

# --- paragraph 24 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Depth of analysis</w:t></w:r></w:p>')

# --- paragraph 25 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">(200 Words) Conclusions: Summarise key findings of research </w:t></w:r></w:p>')

# --- paragraph 26 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Combined methods and what they tell us about sense of home</w:t></w:r></w:p>')

# --- paragraph 27 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">Reflect on report methodology </w:t></w:r></w:p>')

# --- paragraph 28 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Consider future approaches to ethnographic methods</w:t></w:r></w:p>')

# --- paragraph 29 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Relate key arguments to wider geographical debates</w:t></w:r></w:p>')

# --- paragraph 30 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Reference List</w:t></w:r></w:p>')

# --- paragraph 31 ---
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Appendix: Typed Transcript, other photos, fieldnotes</w:t></w:r></w:p>')

